$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = "41.869.55"
$ws.Range("E2").Value = "  +1.19%  "

# Row 3
$ws.Range("D3").Value = "2.206.64"
$ws.Range("E3").Value = "  +0.71%  "

# Row 4
$ws.Range("E4").Value = "  +0.05%  "

# Row 5
$c = $ws.Range("D5")
$c.NumberFormat = "@"
$c.Value = "250.85"
$c.Style = "Normal"
$ws.Range("E5").Value = "  -1.68%  "

# Row 6
$c = $ws.Range("D6")
$c.NumberFormat = "@"
$c.Value = "0.621"
$c.Style = "Normal"
$ws.Range("E6").Value = "  -0.63%  "

# Row 7
$c = $ws.Range("D7")
$c.NumberFormat = "@"
$c.Value = "67.35"
$c.Style = "Normal"
$ws.Range("E7").Value = "  -0.93%  "

# Row 8
$ws.Range("E8").Value = "  -0.05%  "

# Row 9
$c = $ws.Range("D9")
$c.NumberFormat = "@"
$c.Value = "0.619"
$c.Style = "Normal"
$ws.Range("E9").Value = "  +6.82%  "

# Row 10
$c = $ws.Range("D10")
$c.NumberFormat = "@"
$c.Value = "38.55"
$c.Style = "Normal"
$ws.Range("E10").Value = "  +1.76%  "

# Row 11
$c = $ws.Range("D11")
$c.NumberFormat = "@"
$c.Value = "59.27"
$c.Style = "Normal"
$ws.Range("E11").Value = "  +1.98%  "

# Row 12
$c = $ws.Range("D12")
$c.NumberFormat = "@"
$c.Value = "0.0934"
$c.Style = "Normal"
$ws.Range("E12").Value = "  -0.89%  "

# Row 13
$c = $ws.Range("D13")
$c.NumberFormat = "@"
$c.Value = "7.01"
$c.Style = "Normal"
$ws.Range("E13").Value = "  -1.16%  "

# Row 14
$c = $ws.Range("D14")
$c.NumberFormat = "@"
$c.Value = "0.104"
$c.Style = "Normal"
$ws.Range("E14").Value = "  +0.05%  "

# Row 15
$ws.Range("D15").Value = "2.537.46"
$ws.Range("E15").Value = "  +0.85%  "

# Row 16
$c = $ws.Range("D16")
$c.NumberFormat = "@"
$c.Value = "0.871"
$c.Style = "Normal"
$ws.Range("E16").Value = "  +0.02%  "

# Row 17
$c = $ws.Range("D17")
$c.NumberFormat = "@"
$c.Value = "14.45"
$c.Style = "Normal"
$ws.Range("E17").Value = "  -0.42%  "

# Row 18
$ws.Range("D18").Value = "2.205.75"
$ws.Range("E18").Value = "  -0.42%  "

# Row 19
$ws.Range("D19").Value = "41.784.40"
$ws.Range("E19").Value = "  +1.29%  "

# Row 20
$ws.Range("D20").Value = "0.0₃0958"
$ws.Range("E20").Value = "  +0.56%  "

# Row 21
$c = $ws.Range("D21")
$c.NumberFormat = "@"
$c.Value = "72.23"
$c.Style = "Normal"
$ws.Range("E21").Value = "  +0.41%  "

# Row 22
$ws.Range("E22").Value = "  -1.73%  "

# Row 23
$c = $ws.Range("D23")
$c.NumberFormat = "@"
$c.Value = "230.74"
$c.Style = "Normal"
$ws.Range("E23").Value = "  -0.70%  "

# Row 24
$c = $ws.Range("D24")
$c.NumberFormat = "@"
$c.Value = "2.02"
$c.Style = "Normal"
$ws.Range("E24").Value = "  -3.16%  "

# Row 25
$c = $ws.Range("D25")
$c.NumberFormat = "@"
$c.Value = "3.88"
$c.Style = "Normal"
$ws.Range("E25").Value = "  +1.95%  "

# Row 26
$ws.Range("E26").Value = "  +0.11%  "

# Row 27
$c = $ws.Range("D27")
$c.NumberFormat = "@"
$c.Value = "11.14"
$c.Style = "Normal"
$ws.Range("E27").Value = "  -7.15%  "

# Row 28
$c = $ws.Range("D28")
$c.NumberFormat = "@"
$c.Value = "2.40"
$c.Style = "Normal"
$ws.Range("E28").Value = "  -4.74%  "

# Row 29
$ws.Range("E29").Value = "  -1.24%  "

# Row 30
$ws.Range("E30").Value = "  -1.56%  "

# Row 31
$c = $ws.Range("D31")
$c.NumberFormat = "@"
$c.Value = "166.48"
$c.Style = "Normal"
$ws.Range("E31").Value = "  -1.60%  "

# Row 32
$c = $ws.Range("D32")
$c.NumberFormat = "@"
$c.Value = "20.37"
$c.Style = "Normal"
$ws.Range("E32").Value = "  -1.21%  "

# Row 33
$ws.Range("E33").Value = "  +0.41%  "

# Row 34
$ws.Range("E34").Value = "  +6.98%  "

# Row 35
$c = $ws.Range("D35")
$c.NumberFormat = "@"
$c.Value = "0.0779"
$c.Style = "Normal"
$ws.Range("E35").Value = "  +7.07%  "

# Row 36
$ws.Range("E36").Value = "  -0.31%  "

# Row 37
$c = $ws.Range("D37")
$c.NumberFormat = "@"
$c.Value = "25.95"
$c.Style = "Normal"
$ws.Range("E37").Value = "  +3.09%  "

# Row 38
$c = $ws.Range("D38")
$c.NumberFormat = "@"
$c.Value = "4.57"
$c.Style = "Normal"
$ws.Range("E38").Value = "  -0.83%  "

# Row 39
$ws.Range("E39").Value = "  +1.78%  "

# Row 40
$ws.Range("E40").Value = "  +2.98%  "

# Row 41
$c = $ws.Range("D41")
$c.NumberFormat = "@"
$c.Value = "2.22"
$c.Style = "Normal"
$ws.Range("E41").Value = "  -0.56%  "

# Row 42
$ws.Range("B42").Value = "FTXToken"
$ws.Range("C42").Value = "https://coinranking.com/coin/NfeOYfNcl+ftxtoken-ftt"
$c = $ws.Range("D42")
$c.NumberFormat = "@"
$c.Value = "5.17"
$c.Style = "Normal"
$ws.Range("E42").Value = "  +5.94%  "

# Row 43
$ws.Range("B43").Value = "THORChain"
$ws.Range("C43").Value = "https://coinranking.com/coin/ybmU-kKU+thorchain-rune"
$c = $ws.Range("D43")
$c.NumberFormat = "@"
$c.Value = "5.61"
$c.Style = "Normal"
$ws.Range("E43").Value = "  -2.53%  "

# Row 44
$c = $ws.Range("D44")
$c.NumberFormat = "@"
$c.Value = "11.93"
$c.Style = "Normal"
$ws.Range("E44").Value = "  -2.68%  "

# Row 45
$c = $ws.Range("D45")
$c.NumberFormat = "@"
$c.Value = "61.30"
$c.Style = "Normal"
$ws.Range("E45").Value = "  -4.75%  "

# Row 46
$c = $ws.Range("D46")
$c.NumberFormat = "@"
$c.Value = "0.196"
$c.Style = "Normal"
$ws.Range("E46").Value = "  -4.40%  "

# Row 47
$c = $ws.Range("D47")
$c.NumberFormat = "@"
$c.Value = "8.56"
$c.Style = "Normal"
$ws.Range("E47").Value = "  -0.85%  "

# Row 48
$c = $ws.Range("D48")
$c.NumberFormat = "@"
$c.Value = "0.0996"
$c.Style = "Normal"
$ws.Range("E48").Value = "  -1.74%  "

# Row 49
$ws.Range("E49").Value = "  -0.20%  "

# Row 50
$ws.Range("E50").Value = "  +0.06%  "

# Row 51
$ws.Range("E51").Value = "  +4.76%  "
